# Reorder the "Periodo Mora" / "Valor Mora" table (rows 16-28, cols E & F)
# into descending period order: the newest period (2210) moves to the top
# (row 16) and the oldest (2110) moves to the bottom (row 28). The pairing
# of period <-> value travels together with each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 28

# Capture the current (period, value) pairs before overwriting anything.
$periods = @()
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += $ws.Cells.Item($r, 5).Value2
    $values += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse row order.
$count = $lastRow - $firstRow + 1
for ($i = 0; $i -lt $count; $i++) {
    $targetRow = $firstRow + $i
    $sourceIndex = $count - 1 - $i
    $ws.Cells.Item($targetRow, 5).Value = $periods[$sourceIndex]
    $ws.Cells.Item($targetRow, 6).Value = $values[$sourceIndex]
}
